# Asistencias sheet: move the "Terminal" column so it sits right after "N°",
# and re-tune a few alignments/widths to match the refreshed sn-schedules layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asistencias")

# --- Move the "Terminal" column (was column I) to right after "N°" (column B) ---
$ws.Columns("I:I").Cut()
$ws.Columns("C:C").Insert()

# --- Row 4 (the sample data row): tidy up alignment for the moved/shifted columns ---
# New "Terminal" data cell (C4) should look like the "N°" cell next to it (centered both ways)
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108

# "DNI" data cell (now D4, shifted right from old C4) becomes left-aligned like the other text columns
$ws.Range("D4").HorizontalAlignment = -4131
$ws.Range("D4").VerticalAlignment = -4108

# "Fecha" data cell (now H4, shifted right from old G4) switches from left- to right-aligned
$ws.Range("H4").HorizontalAlignment = -4152

# "Hora" data cell (J4) switches from left- to right-aligned
$ws.Range("J4").HorizontalAlignment = -4152

# --- Column widths: widen the relocated "Terminal" column, nudge a couple others ---
$ws.Columns("C:C").ColumnWidth = 16.42578125
$ws.Columns("E:E").ColumnWidth = 22.42578125
$ws.Columns("J:J").ColumnWidth = 17.28515625
$ws.Columns("K:K").ColumnWidth = 21.42578125

# --- Selection cursor moves one column left (to B) ---
$ws.Range("B1:B1048576").Select()
